# Consolidate Extracted Excel Files
# Adds new configuration rows to the "Tasks" sheet of Config.xlsx describing
# the new output / consolidation paths used by the workflow, and refreshes
# the Input_PaySlip_Path value to include a trailing backslash.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- New values, written in the exact order that reproduces the target
#     shared-string table layout (new strings are appended as first seen). ---

# Row 16 - Conolidated_Excel_Path (sic - matches source workbook's typo)
$ws.Range("B16").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Output\Consolidated Excel\"

# Row 17 - Processed_Path
$ws.Range("B17").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Output\Processed\"

$ws.Range("A16").Value = "Conolidated_Excel_Path"
$ws.Range("A17").Value = "Processed_Path"

# Row 12 - Input_PaySlip_Path now gets a trailing backslash
$ws.Range("B12").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Input\"

# Row 14 - Output_Excel_Path
$ws.Range("A14").Value = "Output_Excel_Path"
$ws.Range("B14").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Output\"

# Row 15 - Output_Sheet_Name
$ws.Range("A15").Value = "Output_Sheet_Name"
$ws.Range("B15").Value = "Sheet1"

# Row 18 - Update_Excel_Path
$ws.Range("A18").Value = "Update_Excel_Path"
$ws.Range("B18").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Output\Updated Excel\"

# Row 20 - Salary_Prediction_URL
$ws.Range("A20").Value = "Salary_Prediction_URL"
$ws.Range("B20").Value = "https://rpa-demo-api.herokuapp.com/"

# --- Formatting: B16 carries the same cell style used by A11/A12 in the
#     source workbook (a "Normal"-based style with an explicit black font).
#     Best-effort attempt to copy it across. ---
$ws.Range("B16").Style = $ws.Range("A11").Style
$ws.Range("B16").Font.Name = $ws.Range("A11").Font.Name
$ws.Range("B16").Font.Size = $ws.Range("A11").Font.Size
$ws.Range("B16").Font.Bold = $ws.Range("A11").Font.Bold
$ws.Range("B16").Font.Color = $ws.Range("A11").Font.Color

# --- Selection / view state, to match the saved workbook ---
$ws.Range("B20").Select()
